$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SellData")

# Negligible floating-point re-stamp of the prior last row's timestamp
# (source system recomputed it when appending the new record below).
$ws.Range("E18").Value = 45818.70474744213

# Append a new sale record as row 19
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "1AYB-3AYB-5AYB-1L-1AP"
$ws.Range("C19").Value = "1-1-1-1-3"
$ws.Range("D19").Value = 60500
$ws.Range("E19").Value = 45819.82946605967

# Match the date-column style used by the other rows (style index 2 / custom date format)
$ws.Range("E19").Style = $ws.Range("E18").Style
$ws.Range("E19").NumberFormat = $ws.Range("E18").NumberFormat
